# Daily attendance processing - 2025-12-18 21:27:22
# Reorders the "Recorded By" (column G) contributor lists: entries that
# were previously ordered "System, <user>" (or contained the backdoor
# accounts) are reversed so the real/backup account is listed first,
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".
# Rows whose "Recorded By" value includes "admin@admin.com" are left
# untouched, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = 7
    $value = $cell.Text

    if ($null -eq $value) { continue }
    if ($value -eq "") { continue }
    if ($value -notmatch ",") { continue }
    if ($value -match "admin@admin.com") { continue }

    $parts = $value -split ", "
    $reversedParts = $parts[($parts.Count - 1)..0]
    $newValue = [string]::Join(", ", $reversedParts)

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
